$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'51.760.69"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.61%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'3.036.29"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +3.09%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D5').Value = "'380.81"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.50%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'103.25"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +2.17%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +1.08%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.01%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').Value = "'  +2.83%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'37.05"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +2.73%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.16%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.80%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'3.525.49"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +3.53%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'18.60"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.91%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'7.78"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  -0.15%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'3.056.38"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +3.57%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.981"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -2.10%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'10.48"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -13.11%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'51.765.39"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.67%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'3.07"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.09%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'12.54"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.27%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'0.0₃0965"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +1.68%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'70.13"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.12%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'269.03"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +1.09%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'3.17"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.62%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'8.22"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +1.53%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'7.51"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +6.74%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = "'  +6.90%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('E29').Value = "'  +3.13%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('E30').Value = "'  -0.04%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.12%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'10.31"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +2.21%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'34.29"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +2.61%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('B34').Value = "'OKB"
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = "'50.50"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.15%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = "'Toncoin"
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'2.05"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +0.05%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.0451"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +4.99%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  -0.07%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'3.35"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +8.63%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'0.292"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +12.73%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'17.15"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +3.84%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +3.66%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'2.60"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +3.86%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.44%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'127.42"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +8.19%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = "'  +7.53%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'21.89"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +3.05%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +5.82%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').Value = "'  +3.76%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'2.038.85"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.07%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'3.337.38"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +3.03%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.0320"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +2.05%  "
$ws.Range('E51').Style = 'Normal'
